# "update script for AF" - add an "activity_feed" worksheet after "authentication"
# with seeded content/image/attack columns used by the automation script.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# New sheet right after "authentication"; it becomes the active/visible tab.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "activity_feed"

# Header row.
$ws2.Range("A1").Value = "content"
$ws2.Range("B1").Value = "image"
$ws2.Range("C1").Value = "attack"

$bodyText = "Generated by automation script: Lorem ipsum dolor sit amet, consectetur adipiscing elit. Nam consectetur urna quis lacus volutpat, ut ornare nisi vestibulum. Vivamus malesuada porttitor scelerisque. Donec pellentesque cursus mi, id mollis metus tincidunt ut. In eu elementum dui, et commodo mi. Etiam ultrices diam in ante convallis porta. Phasellus vulputate sagittis pulvinar. Donec id velit facilisis, blandit dui nec, tempor velit. Morbi magna ante, condimentum ut diam semper, interdum consectetur urna. Nullam dignissim condimentum viverra. "

# Build the "plain" bigger-font style on row 3 first (placeholder rows use it as-is).
$ws2.Range("A3").Font.Name = "Arial"
$ws2.Range("A3").Font.Size = 14
$ws2.Range("A3").Font.Color = 0
$ws2.Rows.Item(3).RowHeight = 18

# Reuse the same format on the other placeholder rows.
$ws2.Range("A3").Copy()
$ws2.Range("A4").PasteSpecial(-4122)
$ws2.Rows.Item(4).RowHeight = 18
$ws2.Range("A5").PasteSpecial(-4122)
$ws2.Rows.Item(5).RowHeight = 18

# Row 2 holds the long generated text: same font, plus word wrap, taller row.
$ws2.Range("A2").Value = $bodyText
$ws2.Range("A2").PasteSpecial(-4122)
$ws2.Range("A2").WrapText = $true
$ws2.Rows.Item(2).RowHeight = 228

$ws2.Columns.Item(1).ColumnWidth = 54.8

$ws2.Range("A2").Select()
